# Fix Question 1 sheet: D7 corrected value and selection
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Question 1")
$ws2 = $wb.Worksheets.Item("Question 2")

$ws1.Range("D7").Value = 2895046.8971999995

# Question 2: B11 held a stray pasted-in string (the full list of percentages,
# space separated) instead of being split down the column. Replace it with the
# real numeric series for B11:B162.
$values = @(0.51,0.449999999999999,0.2,0.45,0.319999999999999,0.47,0.26,0.44,0.48,0.52,0.32,0.26,0.05,0.26,0.459999999999999,0.29,0.12,0.52,0.449999999999999,0.16,0.15,0.27,0.11,0.18,0.21,0.11,0.19,0.37,0.26,0.25,0.44,0.659999999999999,0.439999999999999,0.38,0.23,0.35000000000000003,0.11,0.48,0.26,0.37,0.41000000000000003,0.52,0.33,0.30000000000000004,0.12,0.11,0.04,0.35000000000000003,0.36000000000000004,0.6000000000000001,0.09000000000000001,0.21000000000000002,0.18,0.05,0.07,0.07,0.57,0.08,0.30000000000000004,0.339999999999999,0.06,0.09,0.37,0.23,0.13,0.09,0.35,0.22,0.68,0.2,0.32,0.19,0.49,0.22,0.47,0.1,0.3,0.16,0.22,0.27,0.28,0.339999999999999,0.11,0.0699999999999999,0.15,0.16,0.26,0.12,0.24,0.37,0.39,0.56,0.329999999999999,0.28,0.25,0.8,0.61,0.28,0.12,0.08,0.19,0.08,0.38,0.3,0.12,0.0699999999999999,0.25,0.319999999999999,0.15,0.24,0.06,0.54,0.07,0.66,0.65,0.36,0.36,0.13,0.12,0.29,0.18,0.21,0.5,0.2,0.67,0.21,0.28,0.41,0.47,0.55,0.52,0.27,0.279999999999999,0.42,0.58,0.74,0.04,0.52,0.08,0.37,0.24,0.48,0.17,0.27,0.54,0.56,0.51,0,0.5,0.37,0.22,0.24)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = 11 + $i
    $ws2.Cells.Item($row, 2).Value = $values[$i]
}

# Restore selections: sheet1 -> D7, sheet2 -> B165. Select on sheet1 first
# (this transiently activates it), then re-activate sheet2 and select there
# last so sheet2 ends up the active/visible tab again, matching the workbook.
$ws1.Range("D7").Select()
$ws2.Activate()
$ws2.Range("B165").Select()
